# Auto-generated edit script: update FFXIV market-price derived columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled market refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M12").Value = -302.5
$ws.Range("I12").Value = 472.5
$ws.Range("H12").Value = 531.6667
$ws.Range("K12").Value = 472.5

$ws.Range("J39").Value = 90
$ws.Range("I39").Value = 85
$ws.Range("N39").Value = -862
$ws.Range("H39").Value = 86.666664
$ws.Range("M39").Value = 41
$ws.Range("K39").Value = 255
$ws.Range("L39").Value = 270

$ws.Range("K41").Value = 171.25
$ws.Range("L41").Value = 586.0769
$ws.Range("J41").Value = 586.0769
$ws.Range("I41").Value = 171.25
$ws.Range("N41").Value = -1466.0769
$ws.Range("H41").Value = 488.47058
$ws.Range("M41").Value = 268.75

$ws.Range("K86").Value = 1499.75
$ws.Range("I86").Value = 1499.75
$ws.Range("H86").Value = 1499.75
$ws.Range("M86").Value = -376.75

$ws.Range("M89").Value = -1882.75
$ws.Range("K89").Value = 7498.75
$ws.Range("I89").Value = 1499.75
$ws.Range("H89").Value = 1499.75

$ws.Range("I116").Value = 27250
$ws.Range("H116").Value = 12962.917
$ws.Range("K116").Value = 27250
$ws.Range("M116").Value = -23808

$ws.Range("I132").Value = 938.0476
$ws.Range("H132").Value = 1104.4131
$ws.Range("M132").Value = -284.1428000000001
$ws.Range("K132").Value = 2814.1428

$ws.Range("I141").Value = 1219329.1
$ws.Range("H141").Value = 968249.3
$ws.Range("M141").Value = -3652807.3
$ws.Range("K141").Value = 3657987.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K2").Value = 2584824.2
$ws.Range("M2").Value = -2584711.2
$ws.Range("L2").Value = 1150
$ws.Range("J2").Value = 1150
$ws.Range("I2").Value = 2584824.2
$ws.Range("N2").Value = -1376
$ws.Range("H2").Value = 1789847.5

$ws.Range("L32").Value = 6095.5
$ws.Range("J32").Value = 6095.5
$ws.Range("N32").Value = -6669.5
$ws.Range("I32").Value = 1824.6351
$ws.Range("H32").Value = 2504.0908
$ws.Range("M32").Value = -1537.6351
$ws.Range("K32").Value = 1824.6351

$ws.Range("L61").Value = 4180.6665
$ws.Range("J61").Value = 4180.6665
$ws.Range("I61").Value = 1142.5
$ws.Range("N61").Value = -4604.6665
$ws.Range("H61").Value = 2024.5483
$ws.Range("M61").Value = -930.5
$ws.Range("K61").Value = 1142.5

$ws.Range("M102").Value = -55.42859999999996
$ws.Range("I102").Value = 1677.4286
$ws.Range("H102").Value = 1971.2778
$ws.Range("K102").Value = 1677.4286

$ws.Range("J116").Value = 1150
$ws.Range("I116").Value = 2584824.2
$ws.Range("N116").Value = -5738
$ws.Range("H116").Value = 1789847.5
$ws.Range("K116").Value = 2584824.2
$ws.Range("M116").Value = -2582530.2
$ws.Range("L116").Value = 1150

$ws.Range("H136").Value = 2024.5483
$ws.Range("M136").Value = -877.5
$ws.Range("K136").Value = 3427.5
$ws.Range("I136").Value = 1142.5
$ws.Range("L136").Value = 12541.9995
$ws.Range("J136").Value = 4180.6665
$ws.Range("N136").Value = -17641.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J3").Value = 1150
$ws.Range("N3").Value = -1378
$ws.Range("I3").Value = 2584824.2
$ws.Range("H3").Value = 1789847.5
$ws.Range("K3").Value = 2584824.2
$ws.Range("M3").Value = -2584710.2
$ws.Range("L3").Value = 1150

$ws.Range("I94").Value = 458.45456
$ws.Range("H94").Value = 636.24
$ws.Range("M94").Value = -7.454560000000015
$ws.Range("K94").Value = 458.45456
$ws.Range("L94").Value = 1940
$ws.Range("J94").Value = 1940
$ws.Range("N94").Value = -2842

$ws.Range("H107").Value = 2225.3845
$ws.Range("K107").Value = 2193
$ws.Range("M107").Value = -273
$ws.Range("L107").Value = 2333.3333
$ws.Range("I107").Value = 2193
$ws.Range("J107").Value = 2333.3333
$ws.Range("N107").Value = -6173.3333

$ws.Range("M134").Value = -9957.785100000001
$ws.Range("K134").Value = 12492.7851
$ws.Range("L134").Value = 4456.9998
$ws.Range("J134").Value = 1485.6666
$ws.Range("I134").Value = 4164.2617
$ws.Range("H134").Value = 3829.4375
$ws.Range("N134").Value = -9526.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 81.25
$ws.Range("N7").Value = -529
$ws.Range("H7").Value = 155.16667
$ws.Range("K7").Value = 81.25
$ws.Range("M7").Value = 31.75
$ws.Range("L7").Value = 303
$ws.Range("J7").Value = 303

$ws.Range("I58").Value = 1403537.9
$ws.Range("H58").Value = 989583.5600000001
$ws.Range("M58").Value = -1403334.9
$ws.Range("K58").Value = 1403537.9

$ws.Range("H107").Value = 417.68182
$ws.Range("K107").Value = 335.82352
$ws.Range("M107").Value = 1584.17648
$ws.Range("L107").Value = 696
$ws.Range("I107").Value = 335.82352
$ws.Range("J107").Value = 696
$ws.Range("N107").Value = -4536

$ws.Range("L122").Value = 19618.2
$ws.Range("J122").Value = 6539.4
$ws.Range("I122").Value = 5318
$ws.Range("N122").Value = -24518.2
$ws.Range("H122").Value = 6190.4287
$ws.Range("M122").Value = -13504
$ws.Range("K122").Value = 15954

$ws.Range("H136").Value = 989583.5600000001
$ws.Range("M136").Value = -4208063.699999999
$ws.Range("K136").Value = 4210613.699999999
$ws.Range("I136").Value = 1403537.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 73.933334
$ws.Range("H4").Value = 73.933334
$ws.Range("M4").Value = -109.800002
$ws.Range("K4").Value = 221.800002
$ws.Range("L4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("L6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("I6").Value = 76.333336
$ws.Range("H6").Value = 76.333336
$ws.Range("M6").Value = -116.000008
$ws.Range("K6").Value = 229.000008
$ws.Range("N6").ClearContents()

$ws.Range("I10").Value = 223
$ws.Range("H10").Value = 223
$ws.Range("K10").Value = 669
$ws.Range("M10").Value = -530

$ws.Range("M12").Value = 133.000001
$ws.Range("L12").Value = 613.2
$ws.Range("J12").Value = 204.4
$ws.Range("N12").Value = -959.2
$ws.Range("I12").Value = 13.333333
$ws.Range("H12").Value = 132.75
$ws.Range("K12").Value = 39.999999

$ws.Range("I32").Value = 849.5
$ws.Range("H32").Value = 887.4286
$ws.Range("M32").Value = -2265.5
$ws.Range("K32").Value = 2548.5

$ws.Range("H107").Value = 700.8261
$ws.Range("L107").Value = 2547
$ws.Range("J107").Value = 849
$ws.Range("N107").Value = -6387

$ws.Range("H130").Value = 2932.5

$ws.Range("L131").Value = 2472.78
$ws.Range("J131").Value = 824.26
$ws.Range("N131").Value = -12552.78
$ws.Range("H131").Value = 824.26

$ws.Range("J132").Value = 2725
$ws.Range("N132").Value = -29585
$ws.Range("H132").Value = 2237.25
$ws.Range("L132").Value = 24525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("H80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("M80").ClearContents()

$ws.Range("L83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("M83").ClearContents()

$ws.Range("M102").Value = -2268.2856
$ws.Range("L102").Value = 2675
$ws.Range("J102").Value = 2675
$ws.Range("N102").Value = -5919
$ws.Range("I102").Value = 3890.2856
$ws.Range("H102").Value = 3448.3635
$ws.Range("K102").Value = 3890.2856

$ws.Range("J132").Value = 3688.7
$ws.Range("N132").Value = -16126.1
$ws.Range("I132").Value = 2961074.5
$ws.Range("H132").Value = 1675254.6
$ws.Range("M132").Value = -8880693.5
$ws.Range("K132").Value = 8883223.5
$ws.Range("L132").Value = 11066.1

$ws.Range("H136").Value = 8990.75
$ws.Range("L136").Value = 26972.25
$ws.Range("J136").Value = 8990.75
$ws.Range("N136").Value = -32072.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L2").Value = 68000
$ws.Range("J2").Value = 68000
$ws.Range("N2").Value = -68224
$ws.Range("H2").Value = 384800

$ws.Range("H46").Value = 2364.4167
$ws.Range("M46").Value = -1152
$ws.Range("K46").Value = 1340
$ws.Range("I46").Value = 1340

$ws.Range("I68").Value = 2792.625
$ws.Range("H68").Value = 3037.889
$ws.Range("M68").Value = -2043.625
$ws.Range("K68").Value = 2792.625

$ws.Range("I71").Value = 2792.625
$ws.Range("H71").Value = 3037.889
$ws.Range("K71").Value = 13963.125
$ws.Range("M71").Value = -10219.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L70").Value = 29039.4
$ws.Range("J70").Value = 29039.4
$ws.Range("H70").Value = 29039.4
$ws.Range("N70").Value = -29669.4

$ws.Range("L73").Value = 29039.4
$ws.Range("J73").Value = 29039.4
$ws.Range("N73").Value = -31223.4
$ws.Range("H73").Value = 29039.4

$ws.Range("L122").Value = 6861.999899999999
$ws.Range("J122").Value = 2287.3333
$ws.Range("N122").Value = -11761.9999
$ws.Range("H122").Value = 79630

$ws.Range("L131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("N131").Value = -35080
$ws.Range("H131").Value = 25000

$ws.Range("H136").Value = 16341860
$ws.Range("M136").Value = -66670206
$ws.Range("K136").Value = 66672756
$ws.Range("I136").Value = 22224252
$ws.Range("L136").Value = 5650.0002
$ws.Range("J136").Value = 1883.3334
$ws.Range("N136").Value = -10750.0002
